# Update Pais sheet: refresh COVID-19 country stats and re-rank a few
# countries whose total-case counts moved relative to their neighbours.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Timestamp banner
$ws.Range("A1").Value2 = "Datos actualizados a 3 de Agosto de 2020 a las 02:11"

# Row 4
$ws.Range("B4").Value2 = 4813640
$ws.Range("C4").Value2 = 49031
$ws.Range("D4").Value2 = 2377861
$ws.Range("E4").Value2 = 2277432
$ws.Range("G4").Value2 = 449
$ws.Range("H4").Value2 = 158347

# Row 5
$ws.Range("E5").Value2 = 755496
$ws.Range("G5").Value2 = 514
$ws.Range("H5").Value2 = 94130

# Row 25
$ws.Range("B25").Value2 = 116884
$ws.Range("C25").Value2 = 285
$ws.Range("D25").Value2 = 101574
$ws.Range("E25").Value2 = 6365

# Row 95
$ws.Range("B95").Value2 = 6855
$ws.Range("C95").Value2 = 62
$ws.Range("E95").Value2 = 1546
$ws.Range("G95").Value2 = 1
$ws.Range("H95").Value2 = 117

# Row 97
$ws.Range("B97").Value2 = 6323
$ws.Range("C97").Value2 = 4
$ws.Range("D97").Value2 = 5115
$ws.Range("E97").Value2 = 1051

# Row 110
$ws.Range("B110").Value2 = 3837
$ws.Range("C110").Value2 = 146
$ws.Range("D110").Value2 = 623
$ws.Range("E110").Value2 = 3131
$ws.Range("G110").Value2 = 3
$ws.Range("H110").Value2 = 83

# Row 114
$ws.Range("A114").Value2 = "Montenegro"
$ws.Range("B114").Value2 = 3258
$ws.Range("C114").Value2 = 60
$ws.Range("D114").Value2 = 1445
$ws.Range("E114").Value2 = 1762
$ws.Range("G114").Value2 = 1
$ws.Range("H114").Value2 = 51

# Row 115
$ws.Range("A115").Value2 = "Somalia"
$ws.Range("B115").Value2 = 3220
$ws.Range("C115").Value2 = 8
$ws.Range("D115").Value2 = 1598
$ws.Range("E115").Value2 = 1529
$ws.Range("H115").Value2 = 93

# Row 116
$ws.Range("A116").Value2 = "Congo"
$ws.Range("B116").Value2 = 3200
$ws.Range("D116").Value2 = 829
$ws.Range("E116").Value2 = 2317
$ws.Range("H116").Value2 = 54

# Row 133
$ws.Range("A133").Value2 = "Surinam"
$ws.Range("B133").Value2 = 1849
$ws.Range("C133").Value2 = 89
$ws.Range("D133").Value2 = 1194
$ws.Range("E133").Value2 = 628
$ws.Range("G133").Value2 = 1
$ws.Range("H133").Value2 = 27

# Row 134
$ws.Range("A134").Value2 = "Sierra Leona"
$ws.Range("B134").Value2 = 1843
$ws.Range("C134").Value2 = 20
$ws.Range("D134").Value2 = 1375
$ws.Range("E134").Value2 = 401
$ws.Range("H134").Value2 = 67

# Row 135
$ws.Range("A135").Value2 = "Benin"
$ws.Range("B135").Value2 = 1805
$ws.Range("D135").Value2 = 1036
$ws.Range("E135").Value2 = 733
$ws.Range("H135").Value2 = 36

# Row 147
$ws.Range("A147").Value2 = "Niger"
$ws.Range("B147").Value2 = 1147
$ws.Range("C147").Value2 = 11
$ws.Range("D147").Value2 = 1032
$ws.Range("E147").Value2 = 46
$ws.Range("H147").Value2 = 69

# Row 148
$ws.Range("A148").Value2 = "Burkina Faso"
$ws.Range("B148").Value2 = 1143
$ws.Range("C148").Value2 = 0
$ws.Range("D148").Value2 = 945
$ws.Range("E148").Value2 = 145
$ws.Range("H148").Value2 = 53

# Row 153
$ws.Range("D153").Value2 = 787
$ws.Range("E153").Value2 = 72

# Row 166
$ws.Range("B166").Value2 = 474
$ws.Range("C166").Value2 = 44
$ws.Range("E166").Value2 = 268
$ws.Range("G166").Value2 = 1
$ws.Range("H166").Value2 = 21

# Row 180
$ws.Range("B180").Value2 = 182
$ws.Range("C180").Value2 = 9
$ws.Range("E180").Value2 = 42

# Row 183
$ws.Range("A183").Value2 = "Barbados"
$ws.Range("B183").Value2 = 132
$ws.Range("C183").Value2 = 10
$ws.Range("D183").Value2 = 98
$ws.Range("E183").Value2 = 27
$ws.Range("H183").Value2 = 7

# Row 184
$ws.Range("A184").Value2 = "San Martin (Parte Holandesa)"
$ws.Range("B184").Value2 = 128
$ws.Range("D184").Value2 = 64
$ws.Range("E184").Value2 = 49
$ws.Range("H184").Value2 = 15

# Row 192
$ws.Range("B192").Value2 = 89
$ws.Range("C192").Value2 = 1
$ws.Range("E192").Value2 = 3
